# Buffer into guest titrations are now performed before host into guest
# titrations: swap the SampleName (B), SamplePrepMethod (C), and
# PipetteConcentration (G) values between each "host into guestNN" row and
# the following "buffer into guestNN" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(6, 7),
    @(8, 9),
    @(10, 11),
    @(12, 13),
    @(14, 15),
    @(16, 17)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in @("B", "C", "G")) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
